$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.818765640258789
$ws.Range("B1").Value = 5.915517807006836
$ws.Range("C1").Value = 6.422956466674805
$ws.Range("D1").Value = 10.04370021820068
$ws.Range("E1").Value = 6.398133277893066
